# Apply the "start paper v1 8/5" revision to the code-documentation sheet.
# The sheet lists Do-files with purpose / notes; this edit reorganizes and
# extends that list with new do-files (oecd_ind_agg.do, consolidate_emp,
# extrap_bilat_activities) while dropping a few stale notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Remove the three stale notes (delete bottom-most row first so the
#        remaining row numbers don't shift under us) ---------------------
$ws.Rows.Item(11).Delete()   # "Germany's monetary values in 2007 ..."
$ws.Rows.Item(7).Delete()    # "Exchange rates adjustments: for countries ..."
$ws.Rows.Item(5).Delete()    # "However, for Slovenia the numbers from ISIC3 ..."

# --- 2. Make room for the new rows (insert bottom-most point first) ------
$ws.Rows.Item(10).Insert()          # new row for consolidate_emp block
$ws.Rows.Item(4).Resize(2).Insert() # two new rows for oecd_ind_agg block

# --- 3. New note for the existing eurostat_ind_agg.do row (row 3) --------
$ws.Cells.Item(3, 4).Value = "Adjust for exchange rates. Mil of EUR/ECU to USD"

# --- 4. New oecd_ind_agg.do block (rows 4-5) ------------------------------
$ws.Cells.Item(4, 2).Value = "oecd_ind_agg.do"
$ws.Cells.Item(4, 3).Value = "combine OECD MNE activities, FDI flows and stocks data"
$ws.Cells.Item(4, 4).Value = "Adjust for exchange rates. Mil of LCU to USD.For countries that haven't joined Euro Zone, use WDI exchange rates; for countries that are in the Euro Zone, use Euro exchange rates after year of adoption. For Euro Zone countries before year of adoption, first use the fixed rate to translate the numbers into LCU and then use the WDI exchange rates to translate into USD"
$ws.Cells.Item(5, 4).Value = "Adjust for outliers: Germany in 2007 (isic3), Slovenia for all years (isic3) 1000 times smaller"

# --- 5. New consolidate_emp block (row 12) --------------------------------
$ws.Cells.Item(12, 2).Value = "consolidate_emp"
$ws.Cells.Item(12, 3).Value = "consolidate two employment variables, n_emp and n_psn_emp"
$ws.Cells.Item(12, 4).Value = "Use n_psn_emp as primary data and supplement with n_emp. Produce corresponding comparison tables for OECD inward/outward, and Eurostat inward."

# --- 6. New extrap_bilat_activities block (rows 17-19, appended) ---------
$ws.Cells.Item(17, 2).Value = "extrap_bilat_activities"
$ws.Cells.Item(17, 4).Value = "drop time series anomalies"
$ws.Cells.Item(18, 4).Value = "supplement missing country in Eurostat with OECD"
$ws.Cells.Item(19, 4).Value = "Filling missing values: (1) use the corresponding variable from the other data source (OECD - ES or ES - OECD)`n(2) impute zero if non-positive FDI stock`n(3) impute zero revenue if zero employment or zero # enterprise (at least one)`n(4) impute zero revenue if all records (emp,ent,rev) of the opposite direction are zero (or missing, at least one zero)`n(5) use `n"

# --- 7. Update the active selection / view to match the edited document --
$ws.Range("D19").Select()
